{"js": "// \"Wireframes version 2.\" -> revert to \"Wireframes version 1.\":\n// the visible text \"Version 2.\" (held in runs \"Versi\"/\"on\"/\" 2\"/\".\") becomes\n// \"Version 1.\" (held in runs \"Version\"/\" 1.\"): the trailing \".\" run is\n// dropped entirely and the spell-checker-split \"Versi\"+\"on\" runs fold back\n// into a single \"Version\" run.\n\nconst body = context.document.body;\n\n// 1) Drop the lone trailing \".\" run.\nconst dot = body.search(\".\", { matchCase: true });\ndot.load(\"items\");\nawait context.sync();\ndot.items[0].insertText(\"\", \"Replace\");\nawait context.sync();\n\n// 2) \" 2\" -> \" 1.\" (rewrites that run's text).\nconst twoRun = body.search(\" 2\", { matchCase: true });\ntwoRun.load(\"items\");\nawait context.sync();\ntwoRun.items[0].insertText(\" 1.\", \"Replace\");\nawait context.sync();\n\n// 3) Re-join the spell-checker-split \"Versi\" + \"on\" into a single run.\n//    Drop \"on\" first, then grow \"Versi\" into \"Version\" -- each step changes\n//    real text content, so Word folds the touched runs back into one.\nconst onRun = body.search(\"on\", { matchCase: true });\nonRun.load(\"items\");\nawait context.sync();\nonRun.items[0].insertText(\"\", \"Replace\");\nawait context.sync();\n\nconst versiRun = body.search(\"Versi\", { matchCase: true });\nversiRun.load(\"items\");\nawait context.sync();\nversiRun.items[0].insertText(\"Version\", \"Replace\");\nawait context.sync();\n", "ps1": "# \"Wireframes version 2.\" -> revert to \"Wireframes version 1.\":\n# the visible text \"Version 2.\" (held in runs \"Versi\"/\"on\"/\" 2\"/\".\") becomes\n# \"Version 1.\" (held in runs \"Version\"/\" 1.\"): the trailing \".\" run is\n# dropped entirely and the spell-checker-split \"Versi\"+\"on\" runs fold back\n# into a single \"Version\" run.\n\n$d = $word.ActiveDocument\n$rng = $d.Content\n\n# 1) Drop the lone trailing \".\" run.\n[void]$rng.Find.Execute(\".\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n\n# 2) \" 2\" -> \" 1.\" (rewrites that run's text).\n$rng.Start = 0\n[void]$rng.Find.Execute(\" 2\", $false, $false, $false, $false, $false, $true, 1, $false, \" 1.\", 2)\n\n# 3) Re-join the spell-checker-split \"Versi\" + \"on\" into a single run.\n#    Drop \"on\" first, then grow \"Versi\" into \"Version\" -- each step changes\n#    real text content, so Word folds the touched runs back into one.\n$rng.Start = 0\n[void]$rng.Find.Execute(\"on\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n$rng.Start = 0\n[void]$rng.Find.Execute(\"Versi\", $false, $false, $false, $false, $false, $true, 1, $false, \"Version\", 2)\n"}
